# Append two new daily rows (2025-11-09 / serial 45970) for both
# stations to the bottom of the day-data log, then leave the selection
# where the user clicked next (H20), matching the recorded session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - 四方坪站 (station "四方坪站", shared-string index 4)
$ws.Range("A18").Value = 45970
$ws.Range("B18").Value = "四方坪站"
$ws.Range("C18").Value = 10063.49
$ws.Range("D18").Value = 8881.64
$ws.Range("E18").Value = 3279.8
$ws.Range("F18").Value = 415

# Row 19 - 高岭站 (station "高岭站", shared-string index 5)
$ws.Range("A19").Value = 45970
$ws.Range("B19").Value = "高岭站"
$ws.Range("C19").Value = 4260.7
$ws.Range("D19").Value = 3757.85
$ws.Range("E19").Value = 1136.31
$ws.Range("F19").Value = 140

# Move the active selection to where the user clicked next, just
# below the newly entered data.
$ws.Range("H20").Select()
